# Repull data, push all data, mean calculation
# Update the dSF column (F) values for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 5
    4  = 4
    6  = -5
    7  = 2
    8  = -6
    9  = 0
    10 = 3
    11 = 1
    12 = 4
    13 = 2
    14 = -1
    15 = 1
    16 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
